# Adds a new "2022-Q3" worksheet (right after "总计") with fund holding
# data, and updates the "总计" (summary) sheet with a new row for 2022-Q3.

$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $addr, $val) {
    # Forces a cell to be written as TEXT even if the value looks numeric
    # (e.g. "005457" or "8.62"), without leaving a stray direct number
    # format behind.
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the structurally
#    identical "2022-Q2" sheet (same headers/styles), then overwrite
#    its data with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

$q3Data = @(
    @(0, "005457", "景顺长城量化小盘股票", "6.57", "93.58", "1.52", "0.0999", 3),
    @(1, "014202", "天弘中证1000指数增强C", "3.69", "94.06", "1.65", "0.0609", 1),
    @(2, "014201", "天弘中证1000指数增强A", "3.68", "94.06", "1.65", "0.0607", 1),
    @(3, "013466", "博时智选量化多因子股票C", "2.28", "92.38", "1.49", "0.0340", 2),
    @(4, "015496", "景顺中证1000指数增强C", "1.83", "92.63", "1.72", "0.0315", 4),
    @(5, "011500", "九泰量化新兴产业混合", "0.58", "93.85", "3.57", "0.0207", 1),
    @(6, "015495", "景顺中证1000指数增强A", "0.69", "92.63", "1.72", "0.0119", 4),
    @(7, "013465", "博时智选量化多因子股票A", "0.49", "92.38", "1.49", "0.0073", 2)
)

$row = 2
foreach ($rec in $q3Data) {
    $q3Sheet.Range("A$row").Value = $rec[0]
    Set-TextCell $q3Sheet "B$row" $rec[1]
    Set-TextCell $q3Sheet "C$row" $rec[2]
    Set-TextCell $q3Sheet "D$row" $rec[3]
    Set-TextCell $q3Sheet "E$row" $rec[4]
    Set-TextCell $q3Sheet "F$row" $rec[5]
    Set-TextCell $q3Sheet "G$row" $rec[6]
    $q3Sheet.Range("H$row").Value = $rec[7]

    if ($row -gt 6) {
        # Rows 7-9 are brand new - copy the formatting of the A column
        # (bordered/bold/centered) from an existing data row so it
        # matches the rest of the table.
        $q3Sheet.Range("A6").Copy()
        $q3Sheet.Range("A$row").PasteSpecial(-4122)
        $q3Sheet.Range("A$row").Value = $rec[0]
    }

    $row++
}

$q3Sheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3
#    and shift the rest down (which happens automatically), then
#    re-normalize the 0-based index column (A) to 0,1,2,3,4.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").Style = "Normal"

$totalData = @(
    @(0, "2022-Q3", 8, 0.33),
    @(1, "2022-Q2", 5, 0.63),
    @(2, "2022-Q1", 5, 1.08),
    @(3, "2021-Q4", 5, 1.15),
    @(4, "2021-Q3", 1, 0.05)
)

$row = 2
foreach ($rec in $totalData) {
    $totalSheet.Range("A$row").Value = $rec[0]
    $totalSheet.Range("B$row").Value = $rec[1]
    $totalSheet.Range("C$row").Value = $rec[2]
    $totalSheet.Range("D$row").Value = $rec[3]
    $row++
}

$totalSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Restore the originally-selected tab (the last sheet, "2021-Q3").
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
